$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "...if the team member requires[BOOKMARK] assistance."
#    -> "...if the team member requires assistance." (bookmark removed
#    from here; it will be re-created further down, after "useful").
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$d.Content.Find.Execute("requires assistance.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "requires assistance.", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Typo fix: "Regular meeting will be help so that" -> "...useful so that"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Regular meeting will be help so that", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Regular meeting will be useful so that", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Re-insert the "_GoBack" bookmark immediately after "useful" in the
#    Communication bullet (a zero-length bookmark, matching Word's usual
#    "last edit location" marker).
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("useful", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ------------------------------------------------------------------
# 4) "...does the work following coding standards" -> "...does the work follow coding standards"
# ------------------------------------------------------------------
$d.Content.Find.Execute("does the work following coding standards", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "does the work follow coding standards", 2) | Out-Null
